# Swap the values of adjacent "~UC_Sets:" header cells that were
# reordered by a recursive re-run of the export tool.
# For each affected worksheet, the cell in column A, row N and the cell
# in column A, row N+1 need to have their text values swapped.

function Swap-Cells($ws, $r1, $r2) {
    $v1 = $ws.Range("A$r1").Value2
    $v2 = $ws.Range("A$r2").Value2
    $ws.Range("A$r1").Value = $v2
    $ws.Range("A$r2").Value = $v1
}

$wb = $excel.ActiveWorkbook

# Cars: rows 1/2 and rows 7/8
$ws = $wb.Worksheets.Item("Cars")
Swap-Cells $ws 1 2
Swap-Cells $ws 7 8

# Cars_2020: rows 1/2
$ws = $wb.Worksheets.Item("Cars_2020")
Swap-Cells $ws 1 2

# CCS+h2: rows 1/2
$ws = $wb.Worksheets.Item("CCS+h2")
Swap-Cells $ws 1 2

# CH_RH: rows 1/2
$ws = $wb.Worksheets.Item("CH_RH")
Swap-Cells $ws 1 2

# IND_fuels: rows 1/2
$ws = $wb.Worksheets.Item("IND_fuels")
Swap-Cells $ws 1 2

# Power_sector: rows 1/2, 10/11, 17/18
$ws = $wb.Worksheets.Item("Power_sector")
Swap-Cells $ws 1 2
Swap-Cells $ws 10 11
Swap-Cells $ws 17 18

# Thermal_gencap: rows 1/2 and 7/8
$ws = $wb.Worksheets.Item("Thermal_gencap")
Swap-Cells $ws 1 2
Swap-Cells $ws 7 8

# TRA_Policy: rows 1/2
$ws = $wb.Worksheets.Item("TRA_Policy")
Swap-Cells $ws 1 2
